$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:D -> B:E
$ws.Columns.Item(1).Insert()

# Header for new column A - copy formatting from the (shifted) B1 header cell
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ID labels for rows 2:32
$ids = @(
    "Hb 19",
    "Hb 20",
    "32 FO1Hepi",
    "33 FO2H",
    "S 17",
    "20 KR1Ph",
    "21 KR2T",
    "22 KR3S",
    "23 KR4V",
    "24 KR5Mt",
    "71 KR1BuTy",
    "72 KR2Sc",
    "73 KR3Mt",
    "74 KR4Fi",
    "75 KR5Mt",
    "Hb 50",
    "Hb 28",
    "Hb 30",
    "KHb 25",
    "26 SR2Fi",
    "Hb 23",
    "Hb 24",
    "Hb 25",
    "Hb 26",
    "Hb 27",
    "KHb 30",
    "KHb 31",
    "27 SH1Ph",
    "28 SH2R",
    "29 SH3V",
    "KS 76"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
